# Update the "RES installed" sheet's installed capacity (column C)
# This represents newly added RES capacity for the IEEE18_2 2025 case.
$wb = $excel.ActiveWorkbook

$wsRes = $wb.Worksheets.Item("RES installed")
$wsRes.Range("C2").Value = 2
$wsRes.Range("C3").Value = 2.5
$wsRes.Range("C4").Value = 1
$wsRes.Range("C5").Value = 1
$wsRes.Range("C6").Value = 1

# Update the selected cell on that sheet to reflect the recorded view state
$wsRes.Range("C5").Select()

# Force a full recalculation so that dependent formulas (Main!B7 and the
# RANDBETWEEN-driven Pg, Winter/Summer S1-S3 sheets that reference
# 'RES installed') pick up the new, non-zero capacities.
$excel.CalculateFullRebuild()
